$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the census-data note row's content entirely (also drops that now-unused
#    shared string); the row itself stays (just blank).
$ws.Range("A2:B2").Clear()

# 2. Swap the "Area" / "(sq. km)" labels between row 4 and row 6, keeping each
#    row's own formatting untouched.
$ws.Range("A4").Value = "(sq. km)"
$ws.Range("A6").Value = "Area"

# 3. Remove the now-empty spacer row (old row 3), shifting rows 4-6 up to 3-5.
$ws.Rows(3).Delete()

# 4. Drop the 1989 / 2002 columns (B and C), shifting the 2014 column (D) into B.
$ws.Range("B4:C4").EntireColumn.Delete()

# 5. All remaining rows now use a taller 20.1pt custom row height.
$ws.Rows("1:5").RowHeight = 20.1

Write-Output "done"
